# Add a "2022-Q1" worksheet between "2021-Q3" and "总计", fill it with the
# new quarter's holdings data, and update the "总计" (summary) sheet with a
# new row for 2022-Q1 (existing 2021-Q3 row shifts down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q3"
# ---------------------------------------------------------------------
$q3Ws = $wb.Worksheets.Item("2021-Q3")
$newWs = $wb.Worksheets.Add($null, $q3Ws)
$newWs.Name = "2022-Q1"

# Header row (row 1), columns B..H
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $newWs.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$hdrRange = $newWs.Range("B1:H1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160
$hdrRange.Borders.LineStyle = 1

# Data rows 2..6
$rows = @(
    @(0, "009225", "天弘中证中美互联网指数（QDII）A", "1.84", "94.90", "5.86", "0.1078", 7),
    @(1, "009226", "天弘中证中美互联网指数（QDII）C", "0.59", "94.90", "5.86", "0.0346", 7),
    @(2, "012751", "建信纳斯达克100指数（QDII）A 美元现汇", "0.34", "88.02", "3.82", "0.0130", 7),
    @(3, "012752", "建信纳斯达克100指数（QDII）C 人民币", "0.34", "88.02", "3.82", "0.0130", 7),
    @(4, "012753", "建信纳斯达克100指数（QDII）C 美元现汇", "0.34", "88.02", "3.82", "0.0130", 7)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $newWs.Cells.Item($r, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Value = $row[0]

    # Columns B, C : plain text (fund code / fund name)
    $bcRange = $newWs.Range($newWs.Cells.Item($r, 2), $newWs.Cells.Item($r, 3))
    $bcRange.NumberFormat = "@"
    $newWs.Cells.Item($r, 2).Value = $row[1]
    $newWs.Cells.Item($r, 3).Value = $row[2]

    # Columns D, E, F, G : numeric-looking text values, keep them as text
    $defgRange = $newWs.Range($newWs.Cells.Item($r, 4), $newWs.Cells.Item($r, 7))
    $defgRange.NumberFormat = "@"
    $newWs.Cells.Item($r, 4).Value = $row[3]
    $newWs.Cells.Item($r, 5).Value = $row[4]
    $newWs.Cells.Item($r, 6).Value = $row[5]
    $newWs.Cells.Item($r, 7).Value = $row[6]

    # Column H : real number
    $newWs.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary worksheet: insert a new row for 2022-Q1
#    above the existing 2021-Q3 row.
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows(2).Insert()

# The existing 2021-Q3 row (now row 3) shifts its "index" value from 0 to 1
$totalWs.Cells.Item(3, 1).Value = 1

$a2 = $totalWs.Cells.Item(2, 1)
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

$b2 = $totalWs.Cells.Item(2, 2)
$b2.NumberFormat = "@"
$b2.Value = "2022-Q1"

$totalWs.Cells.Item(2, 3).Value = 5
$totalWs.Cells.Item(2, 4).Value = 0.18

# ---------------------------------------------------------------------
# 3. Restore the original active sheet ("2021-Q3"), since adding sheets
#    re-focuses the newly created one.
# ---------------------------------------------------------------------
$q3Ws.Activate()
